$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 15.07846107257523
$ws.Range("C2").Value = 8.44894936692949
$ws.Range("D2").Value = 9.482794649331851
$ws.Range("E2").Value = 13.77982124675043
$ws.Range("F2").Value = 31.18773958003544
$ws.Range("J2").Value = 9.838890416809985
$ws.Range("N2").Value = 16.86215702231472
$ws.Range("O2").Value = 23.16503057641341
$ws.Range("B3").Value = 14.54113448007726
$ws.Range("C3").Value = 7.961981145490929
$ws.Range("D3").Value = 9.440132079982583
$ws.Range("E3").Value = 13.72411774595253
$ws.Range("F3").Value = 31.1662710351466
$ws.Range("J3").Value = 9.844916739155243
$ws.Range("N3").Value = 16.91142088730059
$ws.Range("O3").Value = 23.20166491549961
$ws.Range("B4").Value = 14.20295255604336
$ws.Range("C4").Value = 7.647229857235755
$ws.Range("D4").Value = 9.415314629859951
$ws.Range("E4").Value = 13.692666674456
$ws.Range("F4").Value = 31.16229911810463
$ws.Range("J4").Value = 9.850272393522664
$ws.Range("N4").Value = 16.94355136033888
$ws.Range("O4").Value = 23.23064480849447
$ws.Range("B5").Value = 14.06327637199842
$ws.Range("C5").Value = 7.515091603985273
$ws.Range("D5").Value = 9.405554902380945
$ws.Range("E5").Value = 13.68055134426827
$ws.Range("F5").Value = 31.16299647663767
$ws.Range("J5").Value = 9.852871317558765
$ws.Range("N5").Value = 16.95711894110785
$ws.Range("O5").Value = 23.24408008534252
$ws.Range("B6").Value = 14.03997727293985
$ws.Range("C6").Value = 7.49291903440261
$ws.Range("D6").Value = 9.403955864620471
$ws.Range("E6").Value = 13.67858221507888
$ws.Range("F6").Value = 31.1632521126698
$ws.Range("J6").Value = 9.853328022106247
$ws.Range("N6").Value = 16.95940049030864
$ws.Range("O6").Value = 23.24640902508614
$ws.Range("B7").Value = 14.20107608626509
$ws.Range("C7").Value = 7.645463354524271
$ws.Range("D7").Value = 9.41518156591831
$ws.Range("E7").Value = 13.69250043174801
$ws.Range("F7").Value = 31.16229914720028
$ws.Range("J7").Value = 9.85030575719305
$ws.Range("N7").Value = 16.94373241653498
$ws.Range("O7").Value = 23.23081942634349
$ws.Range("B8").Value = 14.89502703548895
$ws.Range("C8").Value = 8.284357798065225
$ws.Range("D8").Value = 9.467803124235692
$ws.Range("E8").Value = 13.76004961356274
$ws.Range("F8").Value = 31.17842633971456
$ws.Range("J8").Value = 9.840624836818359
$ws.Range("N8").Value = 16.87875304272801
$ws.Range("O8").Value = 23.17631289845743
$ws.Range("B9").Value = 16.18197489636041
$ws.Range("C9").Value = 9.409338431093836
$ws.Range("D9").Value = 9.581590364313003
$ws.Range("E9").Value = 13.91388887860089
$ws.Range("F9").Value = 31.28303896781037
$ws.Range("J9").Value = 9.834764589959892
$ws.Range("N9").Value = 16.76622606063953
$ws.Range("O9").Value = 23.1210978949262
$ws.Range("B10").Value = 17.0726613952575
$ws.Range("C10").Value = 10.15497155153088
$ws.Range("D10").Value = 9.671178845861187
$ws.Range("E10").Value = 14.03928598334524
$ws.Range("F10").Value = 31.40414659051223
$ws.Range("J10").Value = 9.838438674156544
$ws.Range("N10").Value = 16.69258312297936
$ws.Range("O10").Value = 23.11226970770469
$ws.Range("B11").Value = 17.46425112024087
$ws.Range("C11").Value = 10.47621611380617
$ws.Range("D11").Value = 9.713128027480735
$ws.Range("E11").Value = 14.09886259454934
$ws.Range("F11").Value = 31.4687541162844
$ws.Range("J11").Value = 9.841835627596641
$ws.Range("N11").Value = 16.6610315156433
$ws.Range("O11").Value = 23.11517968788383
$ws.Range("B12").Value = 17.61046284851115
$ws.Range("C12").Value = 10.59525922326157
$ws.Range("D12").Value = 9.729175382349323
$ws.Range("E12").Value = 14.12177270271668
$ws.Range("F12").Value = 31.49457606352583
$ws.Range("J12").Value = 9.843369229392483
$ws.Range("N12").Value = 16.64936325708241
$ws.Range("O12").Value = 23.11727895438457
$ws.Range("B13").Value = 17.5790675774363
$ws.Range("C13").Value = 10.56973734848652
$ws.Range("D13").Value = 9.725712244858121
$ws.Range("E13").Value = 14.11682328316196
$ws.Range("F13").Value = 31.48895472476961
$ws.Range("J13").Value = 9.843027960018739
$ws.Range("N13").Value = 16.65186379714668
$ws.Range("O13").Value = 23.11678247114969
$ws.Range("B14").Value = 17.47632226750068
$ws.Range("C14").Value = 10.48606219165943
$ws.Range("D14").Value = 9.714445060651016
$ws.Range("E14").Value = 14.10074049870997
$ws.Range("F14").Value = 31.47085139232829
$ws.Range("J14").Value = 9.841956848300313
$ws.Range("N14").Value = 16.66006595974645
$ws.Range("O14").Value = 23.11533240326938
$ws.Range("B15").Value = 17.41311427131692
$ws.Range("C15").Value = 10.43446885049684
$ws.Range("D15").Value = 9.707564403783312
$ws.Range("E15").Value = 14.09093443517527
$ws.Range("F15").Value = 31.45993884747243
$ws.Range("J15").Value = 9.841332931752449
$ws.Range("N15").Value = 16.6651264224239
$ws.Range("O15").Value = 23.11457410116715
$ws.Range("B16").Value = 17.04678597651355
$ws.Range("C16").Value = 10.13361407534816
$ws.Range("D16").Value = 9.668460588273357
$ws.Range("E16").Value = 14.03544225489206
$ws.Range("F16").Value = 31.40011484482707
$ws.Range("J16").Value = 9.838251325911784
$ws.Range("N16").Value = 16.69468426003273
$ws.Range("O16").Value = 23.11221903929683
$ws.Range("B17").Value = 16.81848428242031
$ws.Range("C17").Value = 9.944433091413348
$ws.Range("D17").Value = 9.644770882132606
$ws.Range("E17").Value = 14.00203821258964
$ws.Range("F17").Value = 31.36584407326659
$ws.Range("J17").Value = 9.836802239698718
$ws.Range("N17").Value = 16.71331576683323
$ws.Range("O17").Value = 23.1125494466825
$ws.Range("B18").Value = 16.68589841379105
$ws.Range("C18").Value = 9.833932860037081
$ws.Range("D18").Value = 9.631258261508332
$ws.Range("E18").Value = 13.9830644010614
$ws.Range("F18").Value = 31.34702887183796
$ws.Range("J18").Value = 9.836131247031188
$ws.Range("N18").Value = 16.72421559409184
$ws.Range("O18").Value = 23.1133913657481
$ws.Range("B19").Value = 16.64079245322498
$ws.Range("C19").Value = 9.796230259898625
$ws.Range("D19").Value = 9.626702832438697
$ws.Range("E19").Value = 13.9766817274431
$ws.Range("F19").Value = 31.34081265722795
$ws.Range("J19").Value = 9.835931991857258
$ws.Range("N19").Value = 16.7279376195763
$ws.Range("O19").Value = 23.11378832951455
$ws.Range("B20").Value = 16.8429200388182
$ws.Range("C20").Value = 9.964746663557554
$ws.Range("D20").Value = 9.647281058889384
$ws.Range("E20").Value = 14.00556946520535
$ws.Range("F20").Value = 31.36939955591627
$ws.Range("J20").Value = 9.836939686961472
$ws.Range("N20").Value = 16.71131342635639
$ws.Range("O20").Value = 23.11244679460344
$ws.Range("B21").Value = 17.50655825298745
$ws.Range("C21").Value = 10.51071048128759
$ws.Range("D21").Value = 9.717750184844096
$ws.Range("E21").Value = 14.10545502890244
$ws.Range("F21").Value = 31.47613206948908
$ws.Range("J21").Value = 9.842264757301965
$ws.Range("N21").Value = 16.65764920081985
$ws.Range("O21").Value = 23.1157312497682
$ws.Range("B22").Value = 17.92814504276761
$ws.Range("C22").Value = 10.85233708317559
$ws.Range("D22").Value = 9.764746058324715
$ws.Range("E22").Value = 14.172767440506
$ws.Range("F22").Value = 31.55378763100743
$ws.Range("J22").Value = 9.847185642317889
$ws.Range("N22").Value = 16.62420628187085
$ws.Range("O22").Value = 23.12369099130143
$ws.Range("B23").Value = 17.70428303281442
$ws.Range("C23").Value = 10.67140107951897
$ws.Range("D23").Value = 9.739580718998933
$ws.Range("E23").Value = 14.13666059856419
$ws.Range("F23").Value = 31.51162299396027
$ws.Range("J23").Value = 9.844427785316659
$ws.Range("N23").Value = 16.64190646726256
$ws.Range("O23").Value = 23.11891059629558
$ws.Range("B24").Value = 16.83187677344097
$ws.Range("C24").Value = 9.955568307847257
$ws.Range("D24").Value = 9.646145874248417
$ws.Range("E24").Value = 14.0039722664719
$ws.Range("F24").Value = 31.36778935671551
$ws.Range("J24").Value = 9.836877042042765
$ws.Range("N24").Value = 16.7122180977928
$ws.Range("O24").Value = 23.1124911728433
$ws.Range("B25").Value = 15.84281578287486
$ws.Range("C25").Value = 9.119040710064523
$ws.Range("D25").Value = 9.549720652922057
$ws.Range("E25").Value = 13.87004790491671
$ws.Range("F25").Value = 31.24694186168523
$ws.Range("J25").Value = 9.834946569094265
$ws.Range("N25").Value = 16.79507811521203
$ws.Range("O25").Value = 23.13047551476544
